# Apply updated "想去人数" (want-to-go count) figures and a refreshed cover
# image URL, matching the regenerated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 0
$ws1.Range("F3").Value  = 22
$ws1.Range("F4").Value  = 16170
$ws1.Range("F6").Value  = 0
$ws1.Range("F7").Value  = 0
$ws1.Range("F8").Value  = 15555
$ws1.Range("F10").Value = 9186
$ws1.Range("F11").Value = 450
$ws1.Range("F12").Value = 0
$ws1.Range("F14").Value = 0
$ws1.Range("F15").Value = 212
$ws1.Range("F17").Value = 0
$ws1.Range("F19").Value = 80
$ws1.Range("F20").Value = 589
$ws1.Range("F24").Value = 0
$ws1.Range("F25").Value = 0
$ws1.Range("F28").Value = 515
$ws1.Range("F30").Value = 0
$ws1.Range("F31").Value = 0
$ws1.Range("F32").Value = 0
$ws1.Range("F33").Value = 0
$ws1.Range("F36").Value = 350
$ws1.Range("F37").Value = 0
$ws1.Range("F38").Value = 0
$ws1.Range("F39").Value = 0
$ws1.Range("F40").Value = 0

$ws1.Range("I40").Value = "//i1.hdslb.com/bfs/openplatform/202409/u3RjLCRL1727662424227.jpeg"

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 0
$ws4.Range("F3").Value  = 22
$ws4.Range("F4").Value  = 16170
$ws4.Range("F5").Value  = 0
$ws4.Range("F6").Value  = 0
$ws4.Range("F8").Value  = 0
$ws4.Range("F10").Value = 9186
$ws4.Range("F11").Value = 450
$ws4.Range("F15").Value = 0
$ws4.Range("F17").Value = 0
$ws4.Range("F18").Value = 0
$ws4.Range("F19").Value = 80
$ws4.Range("F20").Value = 590
$ws4.Range("F22").Value = 0
$ws4.Range("F27").Value = 0
$ws4.Range("F28").Value = 0
$ws4.Range("F29").Value = 0
$ws4.Range("F31").Value = 0
$ws4.Range("F34").Value = 0
$ws4.Range("F35").Value = 0
$ws4.Range("F37").Value = 0
$ws4.Range("F38").Value = 350
$ws4.Range("F39").Value = 471
$ws4.Range("F41").Value = 0

$ws4.Range("I43").Value = "//i1.hdslb.com/bfs/openplatform/202409/u3RjLCRL1727662424227.jpeg"
